# Updates cryptos list D (Price) and E (Volume(1h)) columns for rows 2-51
# D-column values are text that look numeric, so force text via NumberFormat
# "@" before assignment, then reset the style back to Normal so no stray
# cell style ("s" attribute) is left behind, matching the source data shape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.215.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.78%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.572.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.46%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.63%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.595"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.577.71"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.12%  "

$ws.Range("E10").Value = "  -2.86%  "

$ws.Range("E11").Value = "  -1.41%  "

$ws.Range("E12").Value = "  +12.41%  "

$ws.Range("E13").Value = "  +2.73%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.028.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.44%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.239.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.37%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000136"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.30%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.577.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.94%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "335.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.13%  "

$ws.Range("E25").Value = "  +6.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.159"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.28%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0771"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.64%  "

$ws.Range("E30").Value = "  +0.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.26%  "

$ws.Range("E32").Value = "  -3.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "157.65"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.22%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.75%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.899"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("E37").Value = "  -0.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.64%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.852"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.74%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "290.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "135.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.98%  "

$ws.Range("E44").Value = "  +0.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0973"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.74%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.590"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.31%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0530"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.78%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0233"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.91%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.58%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.943.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.79%  "
